$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (Mauvaises herbes / Plante nocive): shorten the "Cueillette manuelle ou..." solution
# text, removing the "Faites bien attention à respecter ce dosage pour ne pas acidifier votre sol." clause.
$ws.Range("C16").Value = "Cueillette manuelle ou Ajoutez dans un litre d'eau, 250 grammes de bicarbonate de soude et 50 ml de vinaigre blanc. Mélangez le tout dans un pulvérisateur, et arrosez les adventices de cette potion."

# Row 15 (Mildiou / Maladies cryptogamiques): shorten the "Diminuer..." solution text,
# removing the "(idéalement de l'eau de pluie ou de l'eau minérale qui est moins calcaire)" clause.
$ws.Range("C15").Value = "Diminuer le taux d'humidité général , Diluer 1 cuillère à soupe de bicarbonate dans 1 litre d'eau, Par temps sec, pulvériser cette solution sur toutes les feuilles en prenant soin de bien atteindre le dessus et le dessous des feuilles."

# That edit also dropped the cell's border/center formatting, leaving only top-aligned wrap text.
$ws.Range("C15").Borders.LineStyle = -4142
$ws.Range("C15").HorizontalAlignment = 1
$ws.Range("C15").VerticalAlignment = -4160
$ws.Range("C15").WrapText = $true
